$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: F4 switch (SW1), micro pin unknown
$ws.Range("A11").Value = "F"
$ws.Range("B11").Value = 4
$ws.Range("D11").Value = "?"
$ws.Range("E11").Value = "SW1"

# Row 12: F0 switch (SW2), micro pin unknown
$ws.Range("A12").Value = "F"
$ws.Range("B12").Value = 0
$ws.Range("D12").Value = "?"
$ws.Range("E12").Value = "SW2"

# Project use for both new rows ("Not used" added after SW1/SW2 strings)
$ws.Range("G11").Value = "Not used"
$ws.Range("G12").Value = "Not used"

# Update the active selection to match the authored state
$ws.Range("G13").Select()
